$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The "OPLEX-N SYRUP 125ML" shortage-report line (originally row 14) was
# removed. The rows below it (ZESTRIL 20MG 10 TAB, the two Arabic-named
# items, the grand-total line and the generated-on footer) all move up by
# one data row, while each row keeps its own original row height/formatting
# (only the cell contents shift - exactly what the authoring tool produced).
# ---------------------------------------------------------------------------

# STEP 1 - pull the data that used to live one row further down up into
# rows 14-16 (row 17's old item data - "سرنجه دواء" - is handled in step 2
# since that row becomes the new grand-total row). Use copy/paste-values so
# the shared-string typing and numeric formatting of each cell is preserved
# exactly (rather than re-typing literals, which Excel would reinterpret).
# Pasting values un-merges the destination cell as a side effect here, so
# each merged group is re-merged immediately afterwards.
foreach ($col in @("C","H","N","P")) {
    $ws.Range("$col" + "15").Copy() | Out-Null
    $ws.Range("$col" + "14").PasteSpecial(-4163) | Out-Null
}
$ws.Range("C14:G14").Merge()
$ws.Range("H14:K14").Merge()
$ws.Range("N14:O14").Merge()

foreach ($col in @("C","H","L","N","P")) {
    $ws.Range("$col" + "16").Copy() | Out-Null
    $ws.Range("$col" + "15").PasteSpecial(-4163) | Out-Null
}
$ws.Range("C15:G15").Merge()
$ws.Range("H15:K15").Merge()
$ws.Range("L15:M15").Merge()
$ws.Range("N15:O15").Merge()

foreach ($col in @("C","H","L","N","P")) {
    $ws.Range("$col" + "17").Copy() | Out-Null
    $ws.Range("$col" + "16").PasteSpecial(-4163) | Out-Null
}
$ws.Range("C16:G16").Merge()
$ws.Range("H16:K16").Merge()
$ws.Range("L16:M16").Merge()
$ws.Range("N16:O16").Merge()

# STEP 2 - row 17 (formerly the last item row) becomes the grand-total row.
$ws.Range("A17:B17").UnMerge()
$ws.Range("C17:G17").UnMerge()
$ws.Range("H17:K17").UnMerge()
$ws.Range("L17:M17").UnMerge()
$ws.Range("N17:O17").UnMerge()
$ws.Range("A17:O17").Clear()

$ws.Range("P18:Q18").Copy() | Out-Null
$ws.Range("P17:Q17").PasteSpecial(-4122) | Out-Null
$ws.Range("P17").Value = 589.29999999999995
$ws.Range("Q17").ClearContents()
$ws.Range("P17:Q17").Merge()

# STEP 3 - row 18 (formerly the grand-total row) becomes the footer row
# that used to be row 19 ("generated on ..." / "page" / "developed by").
$ws.Range("P18:Q18").UnMerge()
$ws.Range("A18:Q18").Clear()

$ws.Range("A19:Q19").Copy() | Out-Null
$ws.Range("A18:Q18").PasteSpecial(-4122) | Out-Null
$ws.Range("A19:Q19").Copy() | Out-Null
$ws.Range("A18:Q18").PasteSpecial(-4163) | Out-Null
$ws.Rows(18).RowHeight = 16.5

$ws.Range("A19:F19").UnMerge()
$ws.Range("G19:I19").UnMerge()
$ws.Range("K19:Q19").UnMerge()
$ws.Range("A18:F18").Merge()
$ws.Range("G18:I18").Merge()
$ws.Range("K18:Q18").Merge()

# Drop the now-duplicated old row 19.
$ws.Rows(19).Delete()
